$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2529582.08
$ws.Range("C7").Value = -43.06690065580031
$ws.Range("D7").Value = 2573
$ws.Range("E7").Value = 2573
$ws.Range("F7").Value = 983.1255654877575
$ws.Range("G7").Value = 4.794076367714628
